$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 103 (existing rows 103-131 shift down to 108-136)
$ws.Rows("103:107").Insert()

# New rows data: Row => (K Variedad, L Calidad, M Volumen, N PrecioMin, O PrecioMax, P PrecioProm,
#                         Q UnidadComercializacion, R Origen, S Precio$/Kg, T Kg/unidad)
$newRows = @{
    103 = @("Castle Brite", "Especial", 75, 15000, 15000, 15000, "`$/caja 10 kilos", "Provincia de San Felipe de Aconcagua", 1500, 10)
    104 = @("Castle Brite", "Primera", 78, 13000, 13000, 13000, "`$/caja 10 kilos", "Provincia de San Felipe de Aconcagua", 1300, 10)
    105 = @("Castle Brite", "Primera", 56, 16000, 16000, 16000, "`$/caja 15 kilos granel", "Provincia de San Felipe de Aconcagua", 1067, 15)
    106 = @("Castle Brite", "Segunda", 70, 10000, 10000, 10000, "`$/caja 10 kilos", "Provincia de San Felipe de Aconcagua", 1000, 10)
    107 = @("Castle Brite", "Segunda", 67, 13000, 13000, 13000, "`$/caja 15 kilos granel", "Provincia de San Felipe de Aconcagua", 867, 15)
}

foreach ($r in 103..107) {
    $data = $newRows[$r]

    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44889
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103003
    $ws.Cells.Item($r, 10).Value = "Damasco"
    $ws.Cells.Item($r, 11).Value = $data[0]
    $ws.Cells.Item($r, 12).Value = $data[1]
    $ws.Cells.Item($r, 13).Value = $data[2]
    $ws.Cells.Item($r, 14).Value = $data[3]
    $ws.Cells.Item($r, 15).Value = $data[4]
    $ws.Cells.Item($r, 16).Value = $data[5]
    $ws.Cells.Item($r, 17).Value = $data[6]
    $ws.Cells.Item($r, 18).Value = $data[7]
    $ws.Cells.Item($r, 19).Value = $data[8]
    $ws.Cells.Item($r, 20).Value = $data[9]
}
